$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the email address (A2) and mobile number (E2) shared-string values.
$ws.Range("A2").Value = "jamefij283a0038@wiroute.com"
$ws.Range("E2").Value = "01122233338"
